$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45898, 5341.91604070659, 3722.14123616529, 1944, 5689.088851, 88.554751935779),
    @(45899, 975.205470885634, 1637.44583015662, 1944, 1881.057642, 24.970750052958),
    @(45900, 862.92289373491, 1584.37896454911, 1944, 1758.686253, 22.3392634922585),
    @(45901, 5594.48067261627, 4901.51975044391, 2952, 6203.576611, 106.608987034485),
    @(45902, 5601.17759246688, 5123.91245185431, 2952, 6210.033796, 115.865360641143),
    @(45903, 5784.11073458408, 5179.80715491303, 2952, 6398.758058, 118.435603263706),
    @(45904, 5812.79253171584, 5162.74419105184, 2952, 6445.449956, 118.475067305667),
    @(45905, 5812.79253171584, 4453.46433942737, 2952, 6445.449956, 88.9217401546473),
    @(45906, 1096.28225605742, 2190.25753228949, 2952, 2113.11063, 10.628579426336),
    @(45907, 970.493609902267, 2169.87892788316, 2952, 1974.370109, 9.23980945753715),
    @(45908, 6110.60514061168, 5086.57515952549, 2952, 6751.167337, 115.630723163076),
    @(45909, 6110.60514061168, 5282.38761243468, 2952, 6751.167337, 123.789575367625),
    @(45910, 6110.60514061168, 5308.30391442236, 2952, 6751.167337, 124.869421283778),
    @(45911, 6110.60514061168, 5326.61005124075, 2952, 6751.167337, 125.632176984545),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}

Write-Output "Done updating need_to_buy data rows 2-15"
